$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values could otherwise be
# auto-converted to numbers by Excel (losing the original text formatting,
# e.g. trailing zeros / exact decimal representation).

$ws.Range("D2").Value = "67.225.19"
$ws.Range("E2").Value = "  -1.55%  "

$ws.Range("D3").Value = "2.486.95"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.58"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.47"
$ws.Range("E6").Value = "  -6.67%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -3.01%  "

$ws.Range("D9").Value = "2.486.28"
$ws.Range("E9").Value = "  -1.92%  "

$ws.Range("E10").Value = "  -3.83%  "

$ws.Range("E11").Value = "  +0.20%  "

$ws.Range("E12").Value = "  -1.83%  "

$ws.Range("E13").Value = "  -3.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.93"
$ws.Range("E14").Value = "  -3.82%  "

$ws.Range("D15").Value = "2.938.84"
$ws.Range("E15").Value = "  -2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -3.59%  "

$ws.Range("D17").Value = "67.066.98"
$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("D18").Value = "2.498.16"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.61"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.80"
$ws.Range("E20").Value = "  -2.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.91"
$ws.Range("E21").Value = "  -2.55%  "

$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.42"
$ws.Range("E23").Value = "  -6.54%  "

$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.51"
$ws.Range("E25").Value = "  -0.36%  "

$ws.Range("E26").Value = "  -4.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.42"
$ws.Range("E27").Value = "  -8.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  -6.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.04"
$ws.Range("E31").Value = "  -3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "500.30"
$ws.Range("E32").Value = "  -7.64%  "

$ws.Range("E33").Value = "  -2.34%  "

$ws.Range("E34").Value = "  -5.37%  "

$ws.Range("E35").Value = "  -0.02%  "

$ws.Range("E36").Value = "  -2.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.36"
$ws.Range("E37").Value = "  +1.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.05"
$ws.Range("E38").Value = "  +0.89%  "

$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("E41").Value = "  -4.49%  "

$ws.Range("E42").Value = "  -5.48%  "

$ws.Range("E43").Value = "  -5.83%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("E45").Value = "  -4.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.35"
$ws.Range("E46").Value = "  -1.52%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.75"
$ws.Range("E47").Value = "  -4.02%  "

$ws.Range("E48").Value = "  -2.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.538"
$ws.Range("E49").Value = "  -4.44%  "

$ws.Range("E50").Value = "  -5.56%  "

$ws.Range("E51").Value = "  -3.90%  "
